# Updates the cryptos list data (Price and Volume(1h) columns) to reflect
# the latest scraped values, mirroring the scheduled GitHub Actions run.
#
# The Price/Volume text is written using a Text-formatted write so that
# numeric-looking strings (e.g. "207.10") are preserved verbatim as text
# rather than being coerced into floating point numbers. The cell style
# is then reset to match an untouched, unstyled data cell so no spurious
# formatting change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cleanStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.699.21"
$ws.Range("D2").Style = $cleanStyle
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("E2").Style = $cleanStyle
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.585.01"
$ws.Range("D3").Style = $cleanStyle
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E3").Style = $cleanStyle
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.07%  "
$ws.Range("E4").Style = $cleanStyle
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.10"
$ws.Range("D5").Style = $cleanStyle
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("E5").Style = $cleanStyle
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E6").Style = $cleanStyle
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E7").Style = $cleanStyle
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.26"
$ws.Range("D8").Style = $cleanStyle
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.98%  "
$ws.Range("E8").Style = $cleanStyle
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E9").Style = $cleanStyle
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("E10").Style = $cleanStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0869"
$ws.Range("D11").Style = $cleanStyle
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E11").Style = $cleanStyle
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.809.23"
$ws.Range("D12").Style = $cleanStyle
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E12").Style = $cleanStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.601.60"
$ws.Range("D13").Style = $cleanStyle
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("E13").Style = $cleanStyle
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("E14").Style = $cleanStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.532"
$ws.Range("D15").Style = $cleanStyle
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.61%  "
$ws.Range("E15").Style = $cleanStyle
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.657.18"
$ws.Range("D16").Style = $cleanStyle
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E16").Style = $cleanStyle
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.25"
$ws.Range("D17").Style = $cleanStyle
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("E17").Style = $cleanStyle
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.11"
$ws.Range("D18").Style = $cleanStyle
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("E18").Style = $cleanStyle
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0693"
$ws.Range("D19").Style = $cleanStyle
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("E19").Style = $cleanStyle
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.52%  "
$ws.Range("E20").Style = $cleanStyle
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E21").Style = $cleanStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("D22").Style = $cleanStyle
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.87%  "
$ws.Range("E22").Style = $cleanStyle
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.63%  "
$ws.Range("E23").Style = $cleanStyle
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("D24").Style = $cleanStyle
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.56%  "
$ws.Range("E24").Style = $cleanStyle
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.57"
$ws.Range("D25").Style = $cleanStyle
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E25").Style = $cleanStyle
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.83"
$ws.Range("D26").Style = $cleanStyle
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E26").Style = $cleanStyle
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("E27").Style = $cleanStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.13"
$ws.Range("D28").Style = $cleanStyle
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("E28").Style = $cleanStyle
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.28%  "
$ws.Range("E29").Style = $cleanStyle
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("E30").Style = $cleanStyle
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0465"
$ws.Range("D31").Style = $cleanStyle
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("E31").Style = $cleanStyle
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("E32").Style = $cleanStyle
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.381.95"
$ws.Range("D33").Style = $cleanStyle
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E33").Style = $cleanStyle
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.06%  "
$ws.Range("E34").Style = $cleanStyle
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("E35").Style = $cleanStyle
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.968"
$ws.Range("D36").Style = $cleanStyle
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("E36").Style = $cleanStyle
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E37").Style = $cleanStyle
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("E38").Style = $cleanStyle
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.538"
$ws.Range("D39").Style = $cleanStyle
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("E39").Style = $cleanStyle
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("E40").Style = $cleanStyle
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("E41").Style = $cleanStyle
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("E42").Style = $cleanStyle
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("E43").Style = $cleanStyle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.53"
$ws.Range("D44").Style = $cleanStyle
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.09%  "
$ws.Range("E44").Style = $cleanStyle
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.03%  "
$ws.Range("E45").Style = $cleanStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.22"
$ws.Range("D46").Style = $cleanStyle
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("E46").Style = $cleanStyle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.720.73"
$ws.Range("D47").Style = $cleanStyle
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("E47").Style = $cleanStyle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.13"
$ws.Range("D48").Style = $cleanStyle
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("E48").Style = $cleanStyle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0998"
$ws.Range("D49").Style = $cleanStyle
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.04%  "
$ws.Range("E49").Style = $cleanStyle
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.43%  "
$ws.Range("E50").Style = $cleanStyle
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.62%  "
$ws.Range("E51").Style = $cleanStyle
